$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$t.Columns.Item(5).Delete()
